# Adding custom models for planets and more.
#
# This script:
#   1. Inserts three new reference sheets (Planets, Countries, Cities) before
#      the existing "Buildings" sheet, each backed by an Excel Table and
#      (where relevant) a drop-down data validation pointing at its parent
#      sheet.
#   2. Adds a "CityTwin" column to the existing Buildings table, linking each
#      building to a row in the new Cities sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Planets sheet
# ---------------------------------------------------------------------
$buildings = $wb.Worksheets.Item("Buildings")

$planets = $wb.Worksheets.Add($buildings)
$planets.Name = "Planets"

$planets.Range("A1").Value = "ID"
$planets.Range("B1").Value = "Name"
$planets.Range("A2").Value = "earth"
$planets.Range("B2").Value = "Earth"
$planets.Range("A3").Value = "mars"
$planets.Range("B3").Value = "Mars"

$planetsTable = $planets.ListObjects.Add(1, $planets.Range("A1:B3"), [System.Reflection.Missing]::Value, 1)
$planetsTable.Name = "Table18"
$planetsTable.TableStyle = "TableStyleLight10"

# ---------------------------------------------------------------------
# 2. Countries sheet
# ---------------------------------------------------------------------
$countries = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $planets)
$countries.Name = "Countries"

$countries.Range("A1").Value = "ID"
$countries.Range("B1").Value = "Planet"
$countries.Range("C1").Value = "Name"

$countries.Range("A2").Value = "netherlands"
$countries.Range("B2").Value = "earth"
$countries.Range("C2").Value = "The Netherlands"

$countries.Range("A3").Value = "germany"
$countries.Range("B3").Value = "earth"
$countries.Range("C3").Value = "Germany"

$countries.Range("A4").Value = "belgium"
$countries.Range("B4").Value = "earth"
$countries.Range("C4").Value = "Belgium"

$countries.Range("A5").Value = "usa"
$countries.Range("B5").Value = "earth"
$countries.Range("C5").Value = "United States Of America"

$countriesTable = $countries.ListObjects.Add(1, $countries.Range("A1:C5"), [System.Reflection.Missing]::Value, 1)
$countriesTable.Name = "Table1811"
$countriesTable.TableStyle = "TableStyleLight10"

$countries.Range("B2:B5").Validation.Add(3, 1, 1, "=Planets!`$A`$2:`$A`$3")
$countries.Range("B2:B5").Validation.InCellDropdown = $true

# ---------------------------------------------------------------------
# 3. Cities sheet
# ---------------------------------------------------------------------
$cities = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $countries)
$cities.Name = "Cities"

$cities.Range("A1").Value = "ID"
$cities.Range("B1").Value = "Country"
$cities.Range("C1").Value = "Name"
$cities.Range("D1").Value = "Population"

$cities.Range("A2").Value = "hilversum"
$cities.Range("B2").Value = "netherlands"
$cities.Range("C2").Value = "Hilversum"
$cities.Range("D2").Value = 90261

$cities.Range("A3").Value = "amsterdam"
$cities.Range("B3").Value = "netherlands"
$cities.Range("C3").Value = "Amsterdam"
$cities.Range("D3").Value = 821752

$cities.Range("A4").Value = "orlando"
$cities.Range("B4").Value = "usa"
$cities.Range("C4").Value = "Orlando"
$cities.Range("D4").Value = 284817

$cities.Range("A5").Value = "antwerp"
$cities.Range("B5").Value = "belgium"
$cities.Range("C5").Value = "Antwerp"
$cities.Range("D5").Value = 506922

$cities.Range("A6").Value = "frankfurt"
$cities.Range("B6").Value = "germany"
$cities.Range("C6").Value = "Frankfurt"
$cities.Range("D6").Value = 753056

$citiesTable = $cities.ListObjects.Add(1, $cities.Range("A1:D6"), [System.Reflection.Missing]::Value, 1)
$citiesTable.Name = "Table181112"
$citiesTable.TableStyle = "TableStyleLight10"

$cities.Range("B2:B6").Validation.Add(3, 1, 1, "=Countries!`$A`$2:`$A`$5")
$cities.Range("B2:B6").Validation.InCellDropdown = $true

# ---------------------------------------------------------------------
# 4. Buildings: insert a "CityTwin" column (after "Name", before "Country")
#    and link each building to a Cities row.
# ---------------------------------------------------------------------
$buildingsTable = $buildings.ListObjects.Item("Table1")
$buildingsTable.ListColumns.Add(3) | Out-Null
$buildings.Range("C1").Value = "CityTwin"
$buildings.Range("C2").Value = "hilversum"
$buildings.Range("C3").Value = "amsterdam"

$buildings.Range("C2:C3").Validation.Add(3, 1, 1, "=Cities!`$A`$2:`$A`$6")
$buildings.Range("C2:C3").Validation.InCellDropdown = $true

# ---------------------------------------------------------------------
# 5. Leave "Buildings" as the active tab, matching the saved file.
# ---------------------------------------------------------------------
$buildings.Activate()
